# "Fruta / hortaliza, semanal"
# A new weekly record is inserted at row 60 (date 2021-10-15 -> serial 44484),
# pushing the existing rows 60-125 down to 61-126 and growing the used range
# from A1:R125 to A1:R126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 60; this shifts rows 60:125
# down to 61:126, carrying all of their values/formatting along with them.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new weekly observation.
$ws.Cells.Item(60, 1).Value  = 5
$ws.Cells.Item(60, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(60, 3).Value  = "Maule"
$ws.Cells.Item(60, 4).Value  = 44484
$ws.Cells.Item(60, 5).Value  = 7
$ws.Cells.Item(60, 6).Value  = 100112017
$ws.Cells.Item(60, 7).Value  = "Apio"
$ws.Cells.Item(60, 8).Value  = "Americana (o)"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 500
$ws.Cells.Item(60, 11).Value = 7000
$ws.Cells.Item(60, 12).Value = 7000
$ws.Cells.Item(60, 13).Value = 7000
$ws.Cells.Item(60, 14).Value = "`$/docena de matas"
$ws.Cells.Item(60, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(60, 16).Value = 1167
$ws.Cells.Item(60, 17).Value = 6
$ws.Cells.Item(60, 18).Value = "Hortaliza"
